$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.05944733333333333
$ws.Range("H2").Value = 0.178342
$ws.Range("M2").Value = 0.05944733333333333
$ws.Range("N2").Value = 0.178342
$ws.Range("Q2").Value = 0.003533985440444444
$ws.Range("R2").Value = 0.031805868964
